$wb = $excel.ActiveWorkbook

# --- Sheet rename / reorder -------------------------------------------------
# Original:  sheetId=12 "xxxTransactionActivity", sheetId=13 "TransactionActivity"
# Target:    sheetId=12 "TransactionActivity",    sheetId=13 "xxTransactionActivity"
# (the physical sheets keep their sheetId / position - only names swap)
$wsNew = $wb.Worksheets("xxxTransactionActivity")
$wsOld = $wb.Worksheets("TransactionActivity")

$wsOld.Name = "xxTransactionActivity"
$wsNew.Name = "TransactionActivity"

# --- Active sheet / selection -----------------------------------------------
# "TransactionActivity" (formerly xxxTransactionActivity) becomes the active /
# selected tab, with D26 selected.
$wsNew.Activate()
$wsNew.Range("D26").Select()

# "xxTransactionActivity" (formerly TransactionActivity) keeps its own
# selection at D16 (unchanged) and is no longer the selected tab - handled
# automatically since only one sheet can be tabSelected at a time.
